$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we are about to update so that Excel
# does not auto-convert numeric-looking / percent-looking strings into
# real numbers (which would lose exact text formatting / precision).
$changedCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "E17", "D18", "E18", "E19", "E20", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "D26", "E26", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "E44", "E45", "E46", "D47", "E47", "D48", "E48", "E49", "E50")
foreach ($addr in $changedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (kept as literal text).
$ws.Range("D2").Value = "261.63"
$ws.Range("E2").Value = "0.91%"
$ws.Range("D3").Value = "27.22"
$ws.Range("E3").Value = "1.06%"
$ws.Range("D4").Value = "4.711"
$ws.Range("E4").Value = "0.21%"
$ws.Range("D5").Value = "0.06204"
$ws.Range("E5").Value = "2.89%"
$ws.Range("D6").Value = "6.717"
$ws.Range("E6").Value = "0.57%"
$ws.Range("D7").Value = "0.8497"
$ws.Range("E7").Value = "-1.06%"
$ws.Range("D8").Value = "0.9155"
$ws.Range("E8").Value = "-1.86%"
$ws.Range("D9").Value = "0.1406"
$ws.Range("E9").Value = "0.95%"
$ws.Range("D10").Value = "0.04613"
$ws.Range("E10").Value = "-5.02%"
$ws.Range("D11").Value = "0.07087"
$ws.Range("E11").Value = "0.30%"
$ws.Range("D12").Value = "0.03151"
$ws.Range("E12").Value = "-0.15%"
$ws.Range("E13").Value = "-0.68%"
$ws.Range("D14").Value = "0.001529"
$ws.Range("E14").Value = "-0.65%"
$ws.Range("D15").Value = "0.0006140"
$ws.Range("E15").Value = "1.31%"
$ws.Range("D16").Value = "0.006044"
$ws.Range("E16").Value = "0.67%"
$ws.Range("E17").Value = "0.13%"
$ws.Range("D18").Value = "3.168"
$ws.Range("E18").Value = "0.14%"
$ws.Range("E19").Value = "0.65%"
$ws.Range("E20").Value = "0.40%"
$ws.Range("E21").Value = "0.86%"
$ws.Range("D22").Value = "4.090"
$ws.Range("E22").Value = "-0.55%"
$ws.Range("D23").Value = "0.04250"
$ws.Range("E23").Value = "0.73%"
$ws.Range("D24").Value = "0.001213"
$ws.Range("E24").Value = "-0.23%"
$ws.Range("E25").Value = "-5.89%"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").Value = "0.19%"
$ws.Range("E40").Value = "1.96%"
$ws.Range("D41").Value = "0.1114"
$ws.Range("E41").Value = "0.00%"
$ws.Range("D42").Value = "0.004138"
$ws.Range("E42").Value = "5.08%"
$ws.Range("D43").Value = "0.002184"
$ws.Range("E43").Value = "-4.72%"
$ws.Range("E44").Value = "-8.83%"
$ws.Range("E45").Value = "1.36%"
$ws.Range("E46").Value = "0.16%"
$ws.Range("D47").Value = "0.03590"
$ws.Range("E47").Value = "-34.17%"
$ws.Range("D48").Value = "0.1677"
$ws.Range("E48").Value = "23.47%"
$ws.Range("E49").Value = "0.16%"
$ws.Range("E50").Value = "0.16%"
